$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The MQTT topic strings used to have a trailing "/A" (an artifact from an
# earlier pin-naming scheme). Strip that suffix everywhere it is used so the
# topics read cleanly (e.g. "fop/decisionRequest/A" -> "fop/decisionRequest").
# Updating column C first (top to bottom) and then column H (top to bottom)
# reproduces the exact order in which the distinct strings were (re)created
# in the shared string table.

$ws.Range("C2").Value = "fop/decisionRequest"
$ws.Range("C3").Value = "fop/summon"
$ws.Range("C4").Value = "fop/resetDecisions"
$ws.Range("C5").Value = "fop/startup"
$ws.Range("C6").Value = "fop/decisionRequest"
$ws.Range("C7").Value = "fop/summon"
$ws.Range("C10").Value = "fop/decisionRequest"
$ws.Range("C11").Value = "fop/summon"
$ws.Range("C12").Value = "fop/resetDecisions"
$ws.Range("C13").Value = "fop/startup"
$ws.Range("C14").Value = "fop/decisionRequest"
$ws.Range("C15").Value = "fop/summon"
$ws.Range("C18").Value = "fop/decisionRequest"
$ws.Range("C19").Value = "fop/summon"
$ws.Range("C20").Value = "fop/resetDecisions"
$ws.Range("C21").Value = "fop/startup"
$ws.Range("C22").Value = "fop/decisionRequest"
$ws.Range("C23").Value = "fop/summon"
$ws.Range("C26").Value = "fop/down"
$ws.Range("C27").Value = "fop/down"

$ws.Range("H8").Value = "refbox/decision"
$ws.Range("H9").Value = "refbox/decision"
$ws.Range("H16").Value = "refbox/decision"
$ws.Range("H17").Value = "refbox/decision"
$ws.Range("H24").Value = "refbox/decision"
$ws.Range("H25").Value = "refbox/decision"

# Reflect the author's last-selected range (whole column H) before saving.
$ws.Columns.Item(8).Select()
